$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly indexing bug-fix: the per-quarter error table was off by one
# quarter for rows 2-11 (Q0..Q9). Recomputed values below.

$data = @(
    @(2,  0.3934793239648547, 0.4037288085692023, 0.2175455044138317, 0.4664177359554756, 0.2626622900771862, 11),
    @(3,  0.3075541611523277, 0.3490093155375321, 0.1614430495769377, 0.4017997630374335, 0.2691182110078417, 13),
    @(4,  0.2355384709971136, 0.2831383568828125, 0.1164317108812595, 0.3412209121394226, 0.2578653762121229, 12),
    @(5,  0.3198529258979658, 0.3395672504785434, 0.1389581556796866, 0.3727709158178606, 0.2030610601727263, 9),
    @(6,  0.3074279757554496, 0.3343017199468005, 0.1343402168544869, 0.3665245105780607, 0.2116761409549216, 9),
    @(7,  0.2926088612260327, 0.3274254659629177, 0.1289431368268439, 0.3590865311131063, 0.2207681817058433, 9),
    @(8,  0.317927339586861,  0.3473839361765009, 0.1472957456731997, 0.3837912787873113, 0.2298271460192152, 8),
    @(9,  0.3052654375636107, 0.3373465939268102, 0.1413070844341813, 0.3759083457894774, 0.23693904963195,   7),
    @(10, 0.3173063146244254, 0.3418763137882747, 0.1463084296434522, 0.3825028491965154, 0.2307148046110977, 7),
    @(11, 0.3606037648954714, 0.379698193710172,  0.1799938953220499, 0.4242568742189688, 0.2448480836729456, 6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
}
